$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = [double]"0.06842864744397358"
$ws.Range("J2").Value = [double]"0.06842864744397358"
$ws.Range("M2").Value = [double]"4.889237666666667"
$ws.Range("N2").Value = [double]"14.667713"
$ws.Range("O2").Value = [double]"0.09529921759032918"
$ws.Range("P2").Value = [double]"0.09529921759032917"
$ws.Range("Q2").Value = [double]"0.3738979315747779"
$ws.Range("R2").Value = [double]"3.365081384173001"
$ws.Range("S2").Value = [double]"0.006521196562175161"
$ws.Range("T2").Value = [double]"0.00652119656217516"
$ws.Range("I3").Value = [double]"0.06842864744397358"
$ws.Range("J3").Value = [double]"0.06842864744397358"
$ws.Range("O3").Value = [double]"0.6664673019309815"
$ws.Range("P3").Value = [double]"0.6664673019309812"
$ws.Range("Q3").Value = [double]"2.614824674903778"
$ws.Range("R3").Value = [double]"23.533422074134"
$ws.Range("S3").Value = [double]"0.04560545603677142"
$ws.Range("T3").Value = [double]"0.04560545603677141"
$ws.Range("I4").Value = [double]"0.06842864744397358"
$ws.Range("J4").Value = [double]"0.06842864744397358"
$ws.Range("K4").Value = [double]"1"
$ws.Range("L4").Value = [double]"0.3333333333333333"
$ws.Range("M4").Value = [double]"0.06428533333333333"
$ws.Range("N4").Value = [double]"0.192856"
$ws.Range("O4").Value = [double]"0.001253026010776221"
$ws.Range("P4").Value = [double]"0.001253026010776221"
$ws.Range("Q4").Value = [double]"0.004916135152888889"
$ws.Range("R4").Value = [double]"0.04424521637600001"
$ws.Range("S4").Value = [double]"8.57428751295347E-05"
$ws.Range("T4").Value = [double]"8.574287512953468E-05"
$ws.Range("I5").Value = [double]"0.06842864744397358"
$ws.Range("J5").Value = [double]"0.06842864744397358"
$ws.Range("M5").Value = [double]"12.052907"
$ws.Range("N5").Value = [double]"36.158721"
$ws.Range("O5").Value = [double]"0.2349308184832226"
$ws.Range("P5").Value = [double]"0.2349308184832226"
$ws.Range("Q5").Value = [double]"0.9217299922823334"
$ws.Range("R5").Value = [double]"8.295569930541001"
$ws.Range("S5").Value = [double]"0.01607599815171259"
$ws.Range("T5").Value = [double]"0.01607599815171259"
$ws.Range("I6").Value = [double]"0.06842864744397358"
$ws.Range("J6").Value = [double]"0.06842864744397358"
$ws.Range("K6").Value = [double]"2"
$ws.Range("L6").Value = [double]"0.6666666666666666"
$ws.Range("M6").Value = [double]"0.1051546666666667"
$ws.Range("N6").Value = [double]"0.315464"
$ws.Range("O6").Value = [double]"0.002049635984690702"
$ws.Range("P6").Value = [double]"0.002049635984690701"
$ws.Range("Q6").Value = [double]"0.008041562927111112"
$ws.Range("R6").Value = [double]"0.07237406634400001"
$ws.Range("S6").Value = [double]"0.0001402538181848816"
$ws.Range("T6").Value = [double]"0.0001402538181848816"
$ws.Range("G7").Value = [double]"0.621785"
$ws.Range("H7").Value = [double]"1.865355"
$ws.Range("I7").Value = [double]"0.5563733034589394"
$ws.Range("J7").Value = [double]"0.5563733034589394"
$ws.Range("M7").Value = [double]"4.889237666666667"
$ws.Range("N7").Value = [double]"14.667713"
$ws.Range("O7").Value = [double]"0.09529921759032918"
$ws.Range("P7").Value = [double]"0.09529921759032917"
$ws.Range("Q7").Value = [double]"3.040054642568334"
$ws.Range("R7").Value = [double]"27.360491783115"
$ws.Range("S7").Value = [double]"0.05302194050778371"
$ws.Range("T7").Value = [double]"0.05302194050778371"
$ws.Range("G8").Value = [double]"0.621785"
$ws.Range("H8").Value = [double]"1.865355"
$ws.Range("I8").Value = [double]"0.5563733034589394"
$ws.Range("J8").Value = [double]"0.5563733034589394"
$ws.Range("O8").Value = [double]"0.6664673019309815"
$ws.Range("P8").Value = [double]"0.6664673019309812"
$ws.Range("Q8").Value = [double]"21.26037407846333"
$ws.Range("R8").Value = [double]"191.34336670617"
$ws.Range("S8").Value = [double]"0.3708046144227066"
$ws.Range("T8").Value = [double]"0.3708046144227065"
$ws.Range("G9").Value = [double]"0.621785"
$ws.Range("H9").Value = [double]"1.865355"
$ws.Range("I9").Value = [double]"0.5563733034589394"
$ws.Range("J9").Value = [double]"0.5563733034589394"
$ws.Range("K9").Value = [double]"1"
$ws.Range("L9").Value = [double]"0.3333333333333333"
$ws.Range("M9").Value = [double]"0.06428533333333333"
$ws.Range("N9").Value = [double]"0.192856"
$ws.Range("O9").Value = [double]"0.001253026010776221"
$ws.Range("P9").Value = [double]"0.001253026010776221"
$ws.Range("Q9").Value = [double]"0.03997165598666667"
$ws.Range("R9").Value = [double]"0.35974490388"
$ws.Range("S9").Value = [double]"0.000697150220935543"
$ws.Range("T9").Value = [double]"0.0006971502209355429"
$ws.Range("G10").Value = [double]"0.621785"
$ws.Range("H10").Value = [double]"1.865355"
$ws.Range("I10").Value = [double]"0.5563733034589394"
$ws.Range("J10").Value = [double]"0.5563733034589394"
$ws.Range("M10").Value = [double]"12.052907"
$ws.Range("N10").Value = [double]"36.158721"
$ws.Range("O10").Value = [double]"0.2349308184832226"
$ws.Range("P10").Value = [double]"0.2349308184832226"
$ws.Range("Q10").Value = [double]"7.494316778995"
$ws.Range("R10").Value = [double]"67.448851010955"
$ws.Range("S10").Value = [double]"0.130709235563823"
$ws.Range("T10").Value = [double]"0.130709235563823"
$ws.Range("G11").Value = [double]"0.621785"
$ws.Range("H11").Value = [double]"1.865355"
$ws.Range("I11").Value = [double]"0.5563733034589394"
$ws.Range("J11").Value = [double]"0.5563733034589394"
$ws.Range("K11").Value = [double]"2"
$ws.Range("L11").Value = [double]"0.6666666666666666"
$ws.Range("M11").Value = [double]"0.1051546666666667"
$ws.Range("N11").Value = [double]"0.315464"
$ws.Range("O11").Value = [double]"0.002049635984690702"
$ws.Range("P11").Value = [double]"0.002049635984690701"
$ws.Range("Q11").Value = [double]"0.06538359441333334"
$ws.Range("R11").Value = [double]"0.58845234972"
$ws.Range("S11").Value = [double]"0.001140362743690682"
$ws.Range("T11").Value = [double]"0.001140362743690682"
$ws.Range("E12").Value = [double]"2"
$ws.Range("F12").Value = [double]"0.6666666666666666"
$ws.Range("G12").Value = [double]"0.4193093333333334"
$ws.Range("H12").Value = [double]"1.257928"
$ws.Range("I12").Value = [double]"0.375198049097087"
$ws.Range("J12").Value = [double]"0.375198049097087"
$ws.Range("M12").Value = [double]"4.889237666666667"
$ws.Range("N12").Value = [double]"14.667713"
$ws.Range("O12").Value = [double]"0.09529921759032918"
$ws.Range("P12").Value = [double]"0.09529921759032917"
$ws.Range("Q12").Value = [double]"2.050102986518223"
$ws.Range("R12").Value = [double]"18.450926878664"
$ws.Range("S12").Value = [double]"0.0357560805203703"
$ws.Range("T12").Value = [double]"0.0357560805203703"
$ws.Range("E13").Value = [double]"2"
$ws.Range("F13").Value = [double]"0.6666666666666666"
$ws.Range("G13").Value = [double]"0.4193093333333334"
$ws.Range("H13").Value = [double]"1.257928"
$ws.Range("I13").Value = [double]"0.375198049097087"
$ws.Range("J13").Value = [double]"0.375198049097087"
$ws.Range("O13").Value = [double]"0.6664673019309815"
$ws.Range("P13").Value = [double]"0.6664673019309812"
$ws.Range("Q13").Value = [double]"14.33722795059022"
$ws.Range("R13").Value = [double]"129.035051555312"
$ws.Range("S13").Value = [double]"0.2500572314715035"
$ws.Range("T13").Value = [double]"0.2500572314715034"
$ws.Range("E14").Value = [double]"2"
$ws.Range("F14").Value = [double]"0.6666666666666666"
$ws.Range("G14").Value = [double]"0.4193093333333334"
$ws.Range("H14").Value = [double]"1.257928"
$ws.Range("I14").Value = [double]"0.375198049097087"
$ws.Range("J14").Value = [double]"0.375198049097087"
$ws.Range("K14").Value = [double]"1"
$ws.Range("L14").Value = [double]"0.3333333333333333"
$ws.Range("M14").Value = [double]"0.06428533333333333"
$ws.Range("N14").Value = [double]"0.192856"
$ws.Range("O14").Value = [double]"0.001253026010776221"
$ws.Range("P14").Value = [double]"0.001253026010776221"
$ws.Range("Q14").Value = [double]"0.02695544026311111"
$ws.Range("R14").Value = [double]"0.242598962368"
$ws.Range("S14").Value = [double]"0.0004701329147111438"
$ws.Range("T14").Value = [double]"0.0004701329147111437"
$ws.Range("E15").Value = [double]"2"
$ws.Range("F15").Value = [double]"0.6666666666666666"
$ws.Range("G15").Value = [double]"0.4193093333333334"
$ws.Range("H15").Value = [double]"1.257928"
$ws.Range("I15").Value = [double]"0.375198049097087"
$ws.Range("J15").Value = [double]"0.375198049097087"
$ws.Range("M15").Value = [double]"12.052907"
$ws.Range("N15").Value = [double]"36.158721"
$ws.Range("O15").Value = [double]"0.2349308184832226"
$ws.Range("P15").Value = [double]"0.2349308184832226"
$ws.Range("Q15").Value = [double]"5.053896398898667"
$ws.Range("R15").Value = [double]"45.48506759008801"
$ws.Range("S15").Value = [double]"0.08814558476768698"
$ws.Range("T15").Value = [double]"0.088145584767687"
$ws.Range("E16").Value = [double]"2"
$ws.Range("F16").Value = [double]"0.6666666666666666"
$ws.Range("G16").Value = [double]"0.4193093333333334"
$ws.Range("H16").Value = [double]"1.257928"
$ws.Range("I16").Value = [double]"0.375198049097087"
$ws.Range("J16").Value = [double]"0.375198049097087"
$ws.Range("K16").Value = [double]"2"
$ws.Range("L16").Value = [double]"0.6666666666666666"
$ws.Range("M16").Value = [double]"0.1051546666666667"
$ws.Range("N16").Value = [double]"0.315464"
$ws.Range("O16").Value = [double]"0.002049635984690702"
$ws.Range("P16").Value = [double]"0.002049635984690701"
$ws.Range("Q16").Value = [double]"0.0440923331768889"
$ws.Range("R16").Value = [double]"0.3968309985920001"
$ws.Range("S16").Value = [double]"0.0001402538181848816"
$ws.Range("T16").Value = [double]"0.0001402538181848816"
